$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, Week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Row 14 (Murder) ---
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0

# --- Row 15 (Rape) ---
$ws.Range("N15").Value = -64.615384615384

# --- Row 16 (Robbery) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 93
$ws.Range("J16").Value = 107
$ws.Range("K16").Value = -13.084112149532
$ws.Range("L16").Value = -10.576923076923
$ws.Range("M16").Value = -66.181818181818
$ws.Range("N16").Value = -89.643652561247

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -3.225806451612
$ws.Range("I17").Value = 312
$ws.Range("J17").Value = 292
$ws.Range("K17").Value = 6.849315068493
$ws.Range("L17").Value = 22.35294117647
$ws.Range("M17").Value = 10.247349823321
$ws.Range("N17").Value = -46.483704974271

# --- Row 18 (Burglary): C/F/H/I/J/K/L/M/N numeric updates ---
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = -13.11475409836
$ws.Range("L18").Value = -10.169491525423
$ws.Range("M18").Value = -82.508250825082
$ws.Range("N18").Value = -94.004524886877
# Row 18: D18/E18 change from text ("0"/"***.*") to numeric (1 / 0).
# Use format-paste from a numeric-styled cell so the cell style (s=14/s=15) updates too,
# then set the numeric value.
$ws.Range("C16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 0

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 141.176470588235
$ws.Range("I19").Value = 354
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = 57.333333333333
$ws.Range("L19").Value = 17.607973421926
$ws.Range("M19").Value = -30.039525691699
$ws.Range("N19").Value = -89.356584485868

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -23.529411764705
$ws.Range("I20").Value = 180
$ws.Range("J20").Value = 167
$ws.Range("K20").Value = 7.784431137724
$ws.Range("L20").Value = 14.649681528662
$ws.Range("M20").Value = -17.43119266055
$ws.Range("N20").Value = -86.97539797395

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -24
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 23.809523809523
$ws.Range("I21").Value = 1021
$ws.Range("J21").Value = 877
$ws.Range("K21").Value = 16.419612314709
$ws.Range("L21").Value = 13.950892857142
$ws.Range("M21").Value = -37.014188772362
$ws.Range("N21").Value = -85.760111576011

# --- Row 23 (Housing) ---
$ws.Range("L23").Value = -33.333333333333

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 59
$ws.Range("H24").Value = 52.542372881355
$ws.Range("I24").Value = 921
$ws.Range("J24").Value = 785
$ws.Range("K24").Value = 17.324840764331
$ws.Range("L24").Value = 15.558343789209
$ws.Range("M24").Value = 10.299401197604

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 273
$ws.Range("J25").Value = 202
$ws.Range("K25").Value = 35.148514851485
$ws.Range("L25").Value = 29.383886255924

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 67
$ws.Range("H26").Value = 28.846153846153
$ws.Range("I26").Value = 559
$ws.Range("J26").Value = 570
$ws.Range("K26").Value = -1.929824561403
$ws.Range("L26").Value = 17.436974789916
$ws.Range("M26").Value = -14.395099540581

# --- Row 28 (Other Sex Crimes): D28/E28 change from numeric to text ---
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 36
$ws.Range("K28").Value = -5.263157894736
$ws.Range("L28").Value = 2.857142857142

# --- Row 29 (Shooting Vic.) ---
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("N29").Value = -92

# --- Row 30 (Shooting Inc.) ---
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("N30").Value = -90.990990990991

# --- Row 33 (Traffic Fatalities): F33 changes from numeric to text "0" ---
$ws.Range("C14").Copy($ws.Range("F33"))
$ws.Range("L33").Value = -42.857142857142

# --- Column width updates (E and H best-fit widths grew) ---
$ws.Columns.Item(5).ColumnWidth = 6.714285714285714
$ws.Columns.Item(8).ColumnWidth = 6.714285714285714

